# edit.ps1
# Adds two new data rows (rows 3 and 4) to the worksheet, replicating the
# result of exploring PSO parameters "a" and "b" separately (one held at
# zero while the other varies), plus data for the "dummy_convergence_plot".
#
# Each new row has 112 columns (A:DH): the same layout as the existing
# row 2 (depth, width, a, b, iw, c1, c2, pop_n, max_iter, conv, dimensions,
# global_best, then 100 "particleN" result columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-separated payload: one line per new row, one field per column (A..DH).
# Each field is prefixed with "N:" (numeric) or "S:" (text) so we know how
# to coerce it when writing into the worksheet.
$tsv = @"
N:3	N:2	N:0.5805749790748964	N:0.7942942102568505	N:0.6116847468948932	N:0.1458389249337071	N:0.3532209859925364	N:10	N:1000	N:0.1	N:7	N:0.8	S:100 , 47.7760 , 4.9053 , 60	S:100 , 130.3691 , 5.5357 , 164	S:100 , 102.5797 , 5.0179 , 129	S:100 , 97.8104 , 5.1732 , 123	S:100 , 122.4535 , 4.3929 , 154	S:100 , 142.2835 , 5.5357 , 179	S:100 , 112.8947 , 5.5357 , 142	S:100 , 142.2835 , 5.5357 , 179	S:100 , 96.2194 , 5.2857 , 121	S:100 , 120.8628 , 4.4643 , 152	S:200 , 50.9471 , 5.1732 , 64	S:200 , 130.3691 , 5.5357 , 164	S:200 , 102.5797 , 5.0179 , 129	S:200 , 93.0531 , 4.8095 , 117	S:200 , 122.4535 , 4.3929 , 154	S:200 , 142.2835 , 5.5357 , 179	S:200 , 130.3691 , 5.5357 , 164	S:200 , 142.2835 , 5.5357 , 179	S:200 , 96.2194 , 5.2857 , 121	S:200 , 112.9198 , 4.4643 , 142	S:300 , 47.7760 , 4.9053 , 60	S:300 , 130.3691 , 5.5357 , 164	S:300 , 102.5797 , 5.0179 , 129	S:300 , 93.0531 , 4.8095 , 117	S:300 , 122.4535 , 4.3929 , 154	S:300 , 142.2835 , 5.5357 , 179	S:300 , 130.3691 , 5.5357 , 164	S:300 , 142.2835 , 5.5357 , 179	S:300 , 99.3966 , 5.2857 , 125	S:300 , 120.8628 , 4.4643 , 152	S:400 , 37.4397 , 5.3815 , 47	S:400 , 130.3691 , 5.5357 , 164	S:400 , 102.5797 , 5.0179 , 129	S:400 , 97.8104 , 5.1732 , 123	S:400 , 122.4535 , 4.3929 , 154	S:400 , 142.2835 , 5.5357 , 179	S:400 , 130.3691 , 5.5357 , 164	S:400 , 142.2835 , 5.5357 , 179	S:400 , 92.2480 , 5.2857 , 116	S:400 , 123.2457 , 4.4643 , 155	S:500 , 50.9471 , 5.1732 , 64	S:500 , 130.3691 , 5.5357 , 164	S:500 , 102.5797 , 5.0179 , 129	S:500 , 93.0531 , 4.8095 , 117	S:500 , 122.4535 , 4.3929 , 154	S:500 , 142.2835 , 5.5357 , 179	S:500 , 130.3691 , 5.5357 , 164	S:500 , 142.2835 , 5.5357 , 179	S:500 , 92.2480 , 5.2857 , 116	S:500 , 120.8578 , 4.6429 , 152	S:600 , 50.9471 , 5.1732 , 64	S:600 , 130.3691 , 5.5357 , 164	S:600 , 102.5797 , 5.0179 , 129	S:600 , 93.0531 , 4.8095 , 117	S:600 , 122.4535 , 4.3929 , 154	S:600 , 142.2835 , 5.5357 , 179	S:600 , 130.3691 , 5.5357 , 164	S:600 , 142.2835 , 5.5357 , 179	S:600 , 96.2194 , 5.2857 , 121	S:600 , 123.2457 , 4.4643 , 155	S:700 , 52.5313 , 5.3815 , 66	S:700 , 130.3691 , 5.5357 , 164	S:700 , 102.5797 , 5.0179 , 129	S:700 , 97.8104 , 5.1732 , 123	S:700 , 122.4535 , 4.3929 , 154	S:700 , 142.2835 , 5.5357 , 179	S:700 , 130.3691 , 5.5357 , 164	S:700 , 142.2835 , 5.5357 , 179	S:700 , 92.2480 , 5.2857 , 116	S:700 , 123.2457 , 4.4643 , 155	S:800 , 45.3640 , 6.5065 , 57	S:800 , 130.3691 , 5.5357 , 164	S:800 , 102.5797 , 5.0179 , 129	S:800 , 93.0514 , 4.8810 , 117	S:800 , 122.4535 , 4.3929 , 154	S:800 , 142.2835 , 5.5357 , 179	S:800 , 130.3691 , 5.5357 , 164	S:800 , 142.2835 , 5.5357 , 179	S:800 , 103.3681 , 5.2857 , 130	S:800 , 112.9198 , 4.4643 , 142	S:900 , 50.9471 , 5.1732 , 64	S:900 , 130.3691 , 5.5357 , 164	S:900 , 102.5797 , 5.0179 , 129	S:900 , 94.6417 , 4.8095 , 119	S:900 , 122.4535 , 4.3929 , 154	S:900 , 142.2835 , 5.5357 , 179	S:900 , 130.3691 , 5.5357 , 164	S:900 , 142.2835 , 5.5357 , 179	S:900 , 92.2480 , 5.2857 , 116	S:900 , 123.2457 , 4.4643 , 155	S:1000 , 50.9471 , 5.1732 , 64	S:1000 , 130.3691 , 5.5357 , 164	S:1000 , 102.5797 , 5.0179 , 129	S:1000 , 97.8104 , 5.1732 , 123	S:1000 , 122.4535 , 4.3929 , 154	S:1000 , 142.2835 , 5.5357 , 179	S:1000 , 130.3691 , 5.5357 , 164	S:1000 , 142.2835 , 5.5357 , 179	S:1000 , 96.2194 , 5.2857 , 121	S:1000 , 120.8628 , 4.4643 , 152
N:3	N:2	N:0.1284803202380342	N:0.04010278080410157	N:0.8432405951299429	N:0.5267426977749526	N:0.7480310273503831	N:10	N:1000	N:0.1	N:7	N:0.8	S:100 , 5.994 , 6.875 , 149	S:100 , 6.961 , 5.529 , 173	S:100 , 6.958 , 6.250 , 173	S:100 , 6.602 , 5.172 , 164	S:100 , 6.602 , 5.172 , 164	S:100 , 6.196 , 6.518 , 154	S:100 , 6.804 , 4.904 , 169	S:100 , 6.961 , 5.529 , 173	S:100 , 6.128 , 3.946 , 152	S:100 , 7.601 , 6.029 , 189	S:200 , 6.398 , 6.029 , 159	S:200 , 6.602 , 5.172 , 164	S:200 , 6.157 , 6.042 , 153	S:200 , 6.725 , 4.696 , 167	S:200 , 6.602 , 5.172 , 164	S:200 , 6.196 , 6.518 , 154	S:200 , 6.725 , 4.696 , 167	S:200 , 6.725 , 4.696 , 167	S:200 , 6.362 , 4.904 , 158	S:200 , 7.601 , 6.029 , 189	S:300 , 6.394 , 7.375 , 159	S:300 , 6.961 , 5.529 , 173	S:300 , 6.958 , 6.250 , 173	S:300 , 6.602 , 5.172 , 164	S:300 , 6.602 , 5.172 , 164	S:300 , 6.196 , 6.518 , 154	S:300 , 6.725 , 4.696 , 167	S:300 , 6.725 , 4.696 , 167	S:300 , 6.721 , 5.417 , 167	S:300 , 7.601 , 6.029 , 189	S:400 , 6.398 , 6.029 , 159	S:400 , 6.602 , 5.172 , 164	S:400 , 6.317 , 6.250 , 157	S:400 , 6.961 , 5.529 , 173	S:400 , 6.961 , 5.529 , 173	S:400 , 5.393 , 6.518 , 134	S:400 , 6.723 , 4.904 , 167	S:400 , 6.723 , 4.904 , 167	S:400 , 6.718 , 6.042 , 167	S:400 , 7.601 , 6.029 , 189	S:500 , 6.474 , 7.375 , 161	S:500 , 6.961 , 5.529 , 173	S:500 , 6.958 , 6.250 , 173	S:500 , 6.602 , 5.172 , 164	S:500 , 6.602 , 5.172 , 164	S:500 , 5.193 , 6.518 , 129	S:500 , 6.523 , 4.904 , 162	S:500 , 6.523 , 4.904 , 162	S:500 , 6.560 , 5.530 , 163	S:500 , 7.601 , 6.029 , 189	S:600 , 6.079 , 5.529 , 151	S:600 , 6.602 , 5.172 , 164	S:600 , 6.958 , 6.250 , 173	S:600 , 6.725 , 4.696 , 167	S:600 , 6.961 , 5.529 , 173	S:600 , 5.394 , 6.310 , 134	S:600 , 6.725 , 4.696 , 167	S:600 , 6.725 , 4.696 , 167	S:600 , 6.562 , 5.017 , 163	S:600 , 7.601 , 6.029 , 189	S:700 , 6.801 , 5.529 , 169	S:700 , 6.602 , 5.172 , 164	S:700 , 6.958 , 6.250 , 173	S:700 , 6.961 , 5.529 , 173	S:700 , 6.961 , 5.529 , 173	S:700 , 5.475 , 6.006 , 136	S:700 , 6.723 , 4.904 , 167	S:700 , 6.723 , 4.904 , 167	S:700 , 6.364 , 4.696 , 158	S:700 , 7.601 , 6.029 , 189	S:800 , 6.074 , 6.875 , 151	S:800 , 6.961 , 5.529 , 173	S:800 , 5.956 , 6.250 , 148	S:800 , 6.961 , 5.529 , 173	S:800 , 6.961 , 5.529 , 173	S:800 , 6.196 , 6.518 , 154	S:800 , 6.523 , 4.904 , 162	S:800 , 6.804 , 4.904 , 169	S:800 , 6.925 , 4.696 , 172	S:800 , 7.601 , 6.029 , 189	S:900 , 6.398 , 6.029 , 159	S:900 , 6.523 , 4.904 , 162	S:900 , 6.958 , 6.250 , 173	S:900 , 6.523 , 4.904 , 162	S:900 , 6.961 , 5.529 , 173	S:900 , 5.193 , 6.518 , 129	S:900 , 6.523 , 4.904 , 162	S:900 , 6.523 , 4.904 , 162	S:900 , 6.720 , 5.529 , 167	S:900 , 7.601 , 6.029 , 189	S:1000 , 6.478 , 6.029 , 161	S:1000 , 6.602 , 5.172 , 164	S:1000 , 6.958 , 6.250 , 173	S:1000 , 6.602 , 5.172 , 164	S:1000 , 6.523 , 4.904 , 162	S:1000 , 5.393 , 6.518 , 134	S:1000 , 6.723 , 4.904 , 167	S:1000 , 6.723 , 4.904 , 167	S:1000 , 5.518 , 5.321 , 137	S:1000 , 7.601 , 6.029 , 189
"@

$rowLines = $tsv -split "`r?`n" | Where-Object { $_.Length -gt 0 }

$startRow = 3
$numCols = 112

$numRows = $rowLines.Length
$data = New-Object 'object[,]' $numRows,$numCols

for ($r = 0; $r -lt $numRows; $r++) {
    $fields = $rowLines[$r] -split "`t"
    for ($c = 0; $c -lt $numCols; $c++) {
        $field = $fields[$c]
        $kind = $field.Substring(0, 2)
        $payload = $field.Substring(2)
        if ($kind -eq "N:") {
            $data[$r, $c] = [double]$payload
        } else {
            $data[$r, $c] = $payload
        }
    }
}

$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$targetRange.Value = $data

# Match the center/center alignment style ("s=1") used by every other
# data cell in the sheet.
$targetRange.VerticalAlignment = -4108
$targetRange.HorizontalAlignment = -4108

Write-Output "Added rows $startRow..$endRow ($numCols columns) to sheet '$($ws.Name)'."
